$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values
$ws.Range("B2").Value = 4.5
$ws.Range("C2").Value = 11
$ws.Range("B3").Value = 4.5
$ws.Range("B4").Value = 0.7

# Add new cell K1 = 1 (extends used range to A1:K5)
$ws.Range("K1").Value = 1
